$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Standard Deviation" header (K4) gets a new neighbour: "Mean" in L4.
# Copy K4's formatting (bold white-on-black header style) onto L4 before
# writing the new header text, so the new column matches the existing ones.
$ws.Range("K4").Copy() | Out-Null
$ws.Range("L4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("L4").Value = "Mean"

# Scroll the sheet right a column and move the active selection from J8 to
# K8, as happens when a user scrolls over to review the new column.
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 2
$ws.Range("K8").Select()
